$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 13839.272
$ws.Range("I32").Value = 18994
$ws.Range("J32").Value = 12693.777
$ws.Range("K32").Value = 18994
$ws.Range("L32").Value = 12693.777
$ws.Range("M32").Value = -18668
$ws.Range("N32").Value = -13345.777
$ws.Range("H112").Value = 3334.2354
$ws.Range("J112").Value = 2883.4062
$ws.Range("L112").Value = 8650.2186
$ws.Range("N112").Value = -10866.2186
$ws.Range("H135").Value = 7322.0835
$ws.Range("I135").Value = 2838.0715
$ws.Range("J135").Value = 13599.7
$ws.Range("K135").Value = 25542.6435
$ws.Range("L135").Value = 122397.3
$ws.Range("M135").Value = -23007.6435
$ws.Range("N135").Value = -127467.3
$ws.Range("H138").Value = 3161.9412
$ws.Range("I138").Value = 1419.8572
$ws.Range("J138").Value = 3733.5625
$ws.Range("K138").Value = 4259.571599999999
$ws.Range("L138").Value = 11200.6875
$ws.Range("M138").Value = 880.4284000000007
$ws.Range("N138").Value = -21480.6875
$ws.Range("H141").Value = 2639.1875
$ws.Range("J141").Value = 3732
$ws.Range("L141").Value = 11196
$ws.Range("N141").Value = -21556

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1570197
$ws.Range("I32").Value = 4298.375
$ws.Range("K32").Value = 4298.375
$ws.Range("M32").Value = -4011.375
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H74").Value = 5997.5713
$ws.Range("I74").Value = 6488.25
$ws.Range("K74").Value = 6488.25
$ws.Range("M74").Value = -5614.25
$ws.Range("H77").Value = 5997.5713
$ws.Range("I77").Value = 6488.25
$ws.Range("K77").Value = 32441.25
$ws.Range("M77").Value = -28073.25
$ws.Range("H132").Value = 640942.3
$ws.Range("I132").Value = 709186.4399999999
$ws.Range("J132").Value = 137641.88
$ws.Range("K132").Value = 2127559.32
$ws.Range("L132").Value = 412925.64
$ws.Range("M132").Value = -2125029.32
$ws.Range("N132").Value = -417985.64

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4431.36
$ws.Range("I86").Value = 2868.6155
$ws.Range("K86").Value = 2868.6155
$ws.Range("M86").Value = -1745.6155
$ws.Range("H89").Value = 4431.36
$ws.Range("I89").Value = 2868.6155
$ws.Range("K89").Value = 14343.0775
$ws.Range("M89").Value = -8727.077499999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8511.652
$ws.Range("I58").Value = 5414.6
$ws.Range("K58").Value = 5414.6
$ws.Range("M58").Value = -5211.6
$ws.Range("H132").Value = 7734.6665
$ws.Range("I132").Value = 4316.0835
$ws.Range("K132").Value = 12948.2505
$ws.Range("M132").Value = -10418.2505
$ws.Range("H136").Value = 8511.652
$ws.Range("I136").Value = 5414.6
$ws.Range("K136").Value = 16243.8
$ws.Range("M136").Value = -13693.8
$ws.Range("H141").Value = 245856.1
$ws.Range("J141").Value = 268173.44
$ws.Range("L141").Value = 268173.44
$ws.Range("N141").Value = -278533.44

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 480.1
$ws.Range("I8").Value = 480.1
$ws.Range("K8").Value = 1440.3
$ws.Range("M8").Value = -1301.3
$ws.Range("H23").Value = 118.82353
$ws.Range("I23").Value = 113.5
$ws.Range("J23").Value = 123.55556
$ws.Range("K23").Value = 340.5
$ws.Range("L23").Value = 370.66668
$ws.Range("M23").Value = -105.5
$ws.Range("N23").Value = -840.66668
$ws.Range("H34").Value = 2460.9443
$ws.Range("I34").Value = 306.7857
$ws.Range("K34").Value = 920.3571000000001
$ws.Range("M34").Value = -836.3571000000001
$ws.Range("H37").Value = 109946.22
$ws.Range("J37").Value = 109946.22
$ws.Range("L37").Value = 329838.66
$ws.Range("N37").Value = -330062.66
$ws.Range("H104").Value = 3511.75
$ws.Range("I104").Value = 424
$ws.Range("J104").Value = 6599.5
$ws.Range("K104").Value = 1272
$ws.Range("L104").Value = 19798.5
$ws.Range("M104").Value = 1349
$ws.Range("N104").Value = -25040.5
$ws.Range("H121").Value = 63286.2
$ws.Range("J121").Value = 101899.336
$ws.Range("L121").Value = 305698.008
$ws.Range("N121").Value = -308318.008
$ws.Range("H122").Value = 2663.4092
$ws.Range("J122").Value = 4043.7856
$ws.Range("L122").Value = 36394.0704
$ws.Range("N122").Value = -41294.0704
$ws.Range("H123").Value = 990
$ws.Range("I123").Value = 990
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 2970
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -520
$ws.Range("N123").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6178.8335
$ws.Range("J80").Value = 10349.2
$ws.Range("L80").Value = 10349.2
$ws.Range("N80").Value = -12345.2
$ws.Range("H83").Value = 6178.8335
$ws.Range("J83").Value = 10349.2
$ws.Range("L83").Value = 51746
$ws.Range("N83").Value = -61730
$ws.Range("H136").Value = 18768.078
$ws.Range("J136").Value = 18768.078
$ws.Range("L136").Value = 56304.234
$ws.Range("N136").Value = -61404.234

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 381666.66
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H41").Value = 33333
$ws.Range("I41").Value = 33333
$ws.Range("K41").Value = 33333
$ws.Range("M41").Value = -32895
$ws.Range("H50").Value = 38832.668
$ws.Range("I50").Value = 39498
$ws.Range("J50").Value = 38500
$ws.Range("K50").Value = 39498
$ws.Range("L50").Value = 38500
$ws.Range("M50").Value = -38861
$ws.Range("N50").Value = -39774
$ws.Range("H55").Value = 3694.8857
$ws.Range("I55").Value = 2691.125
$ws.Range("K55").Value = 2691.125
$ws.Range("M55").Value = -2518.125
$ws.Range("H82").Value = 2585.2856
$ws.Range("I82").Value = 1876.5
$ws.Range("J82").Value = 3335.7646
$ws.Range("K82").Value = 1876.5
$ws.Range("L82").Value = 3335.7646
$ws.Range("M82").Value = -1515.5
$ws.Range("N82").Value = -4057.7646
$ws.Range("H85").Value = 2585.2856
$ws.Range("I85").Value = 1876.5
$ws.Range("J85").Value = 3335.7646
$ws.Range("K85").Value = 1876.5
$ws.Range("L85").Value = 3335.7646
$ws.Range("M85").Value = -628.5
$ws.Range("N85").Value = -5831.7646
$ws.Range("H136").Value = 18314.875
$ws.Range("I136").Value = 28200
$ws.Range("J136").Value = 15019.833
$ws.Range("K136").Value = 84600
$ws.Range("L136").Value = 45059.499
$ws.Range("M136").Value = -82050
$ws.Range("N136").Value = -50159.499

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 41995
$ws.Range("J42").Value = 41995
$ws.Range("L42").Value = 41995
$ws.Range("N42").Value = -42751
$ws.Range("H62").Value = 13290.5
$ws.Range("J62").Value = 11938.5
$ws.Range("L62").Value = 11938.5
$ws.Range("N62").Value = -13186.5
$ws.Range("H65").Value = 13290.5
$ws.Range("J65").Value = 11938.5
$ws.Range("L65").Value = 59692.5
$ws.Range("N65").Value = -65932.5
$ws.Range("H107").Value = 5714832.5
$ws.Range("I107").Value = 6667222.5
$ws.Range("J107").Value = 494.2
$ws.Range("K107").Value = 20001667.5
$ws.Range("L107").Value = 1482.6
$ws.Range("M107").Value = -19999747.5
$ws.Range("N107").Value = -5322.6
$ws.Range("H136").Value = 10879438
$ws.Range("J136").Value = 5755.467
$ws.Range("L136").Value = 17266.401
$ws.Range("N136").Value = -22366.401
